# Auto-generated edit script: updates 想去人数 (F column) counts across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 659
$ws.Range("F3").Value = 6470
$ws.Range("F4").Value = 1036
$ws.Range("F5").Value = 640
$ws.Range("F6").Value = 1433
$ws.Range("F7").Value = 3183
$ws.Range("F9").Value = 547
$ws.Range("F10").Value = 2088
$ws.Range("F11").Value = 449
$ws.Range("F12").Value = 376
$ws.Range("F13").Value = 220
$ws.Range("F15").Value = 233
$ws.Range("F16").Value = 1031
$ws.Range("F17").Value = 403
$ws.Range("F19").Value = 152
$ws.Range("F20").Value = 4023
$ws.Range("F21").Value = 1221
$ws.Range("F22").Value = 3164
$ws.Range("F23").Value = 309
$ws.Range("F24").Value = 96
$ws.Range("F25").Value = 2933
$ws.Range("F26").Value = 2933
$ws.Range("F27").Value = 4590
$ws.Range("F30").Value = 509
$ws.Range("F31").Value = 3015
$ws.Range("F32").Value = 286
$ws.Range("F35").Value = 64
$ws.Range("F37").Value = 1093
$ws.Range("F38").Value = 1341
$ws.Range("F40").Value = 1211
$ws.Range("F41").Value = 783
$ws.Range("F43").Value = 721
$ws.Range("F44").Value = 473
$ws.Range("F46").Value = 193
$ws.Range("F47").Value = 34
$ws.Range("F48").Value = 80
$ws.Range("F49").Value = 341
$ws.Range("F50").Value = 3662

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 957

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1590

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6470
$ws.Range("F4").Value = 640
$ws.Range("F5").Value = 1433
$ws.Range("F6").Value = 3183
$ws.Range("F8").Value = 547
$ws.Range("F10").Value = 2088
$ws.Range("F11").Value = 449
$ws.Range("F12").Value = 376
$ws.Range("F13").Value = 220
$ws.Range("F14").Value = 957
$ws.Range("F17").Value = 233
$ws.Range("F18").Value = 1031
$ws.Range("F20").Value = 403
$ws.Range("F21").Value = 152
$ws.Range("F22").Value = 4023
$ws.Range("F24").Value = 1221
$ws.Range("F26").Value = 3164
$ws.Range("F27").Value = 2933
$ws.Range("F28").Value = 2933
$ws.Range("F29").Value = 4590
$ws.Range("F31").Value = 3015
$ws.Range("F32").Value = 286
$ws.Range("F35").Value = 1093
$ws.Range("F36").Value = 1341
$ws.Range("F38").Value = 1211
$ws.Range("F39").Value = 783
$ws.Range("F41").Value = 473
$ws.Range("F46").Value = 193
$ws.Range("F47").Value = 34
$ws.Range("F48").Value = 80
$ws.Range("F49").Value = 341
$ws.Range("F50").Value = 3662

